$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A1 text and wrap the text for A1:C1
$ws.Range("A1").Value = "Hong Infinity Insulated Glass Inc."

# Apply wrap text alignment to A1:C1 (creates the new cellXfs entry with wrapText)
$ws.Range("A1:C1").WrapText = $true
